$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cD = $ws.Range('D2')
$cD.NumberFormat = '@'
$cD.Value = '51.473.82'
$cD.Style = 'Normal'
$ws.Range('E2').Value = '  -0.03%  '
$cD = $ws.Range('D3')
$cD.NumberFormat = '@'
$cD.Value = '3.051.73'
$cD.Style = 'Normal'
$ws.Range('E3').Value = '  +2.09%  '
$cD = $ws.Range('D4')
$cD.NumberFormat = '@'
$cD.Value = '1.00'
$cD.Style = 'Normal'
$ws.Range('E4').Value = '  +0.14%  '
$cD = $ws.Range('D5')
$cD.NumberFormat = '@'
$cD.Value = '385.20'
$cD.Style = 'Normal'
$ws.Range('E5').Value = '  +0.91%  '
$cD = $ws.Range('D6')
$cD.NumberFormat = '@'
$cD.Value = '103.07'
$cD.Style = 'Normal'
$ws.Range('E6').Value = '  -0.41%  '
$ws.Range('E7').Value = '  -0.54%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  -1.61%  '
$cD = $ws.Range('D10')
$cD.NumberFormat = '@'
$cD.Value = '36.75'
$cD.Style = 'Normal'
$ws.Range('E10').Value = '  -0.01%  '
$ws.Range('E11').Value = '  +0.06%  '
$cD = $ws.Range('D12')
$cD.NumberFormat = '@'
$cD.Value = '0.0861'
$cD.Style = 'Normal'
$ws.Range('E12').Value = '  +0.00%  '
$cD = $ws.Range('D13')
$cD.NumberFormat = '@'
$cD.Value = '3.534.01'
$cD.Style = 'Normal'
$ws.Range('E13').Value = '  +2.12%  '
$ws.Range('E14').Value = '  +0.46%  '
$ws.Range('E15').Value = '  -0.78%  '
$cD = $ws.Range('D16')
$cD.NumberFormat = '@'
$cD.Value = '3.051.38'
$cD.Style = 'Normal'
$ws.Range('E16').Value = '  +1.46%  '
$cD = $ws.Range('D17')
$cD.NumberFormat = '@'
$cD.Value = '0.972'
$cD.Style = 'Normal'
$ws.Range('E17').Value = '  -2.90%  '
$cD = $ws.Range('D18')
$cD.NumberFormat = '@'
$cD.Value = '10.65'
$cD.Style = 'Normal'
$ws.Range('E18').Value = '  -6.00%  '
$cD = $ws.Range('D19')
$cD.NumberFormat = '@'
$cD.Value = '51.570.61'
$cD.Style = 'Normal'
$ws.Range('E19').Value = '  +0.10%  '
$cD = $ws.Range('D20')
$cD.NumberFormat = '@'
$cD.Value = '3.14'
$cD.Style = 'Normal'
$ws.Range('E20').Value = '  -0.21%  '
$cD = $ws.Range('D21')
$cD.NumberFormat = '@'
$cD.Value = '12.42'
$cD.Style = 'Normal'
$ws.Range('E21').Value = '  -1.63%  '
$cD = $ws.Range('D22')
$cD.NumberFormat = '@'
$cD.Value = '0.0₃0965'
$cD.Style = 'Normal'
$ws.Range('E22').Value = '  +0.12%  '
$cD = $ws.Range('D23')
$cD.NumberFormat = '@'
$cD.Value = '70.15'
$cD.Style = 'Normal'
$ws.Range('E23').Value = '  -0.29%  '
$cD = $ws.Range('D24')
$cD.NumberFormat = '@'
$cD.Value = '267.92'
$cD.Style = 'Normal'
$ws.Range('E24').Value = '  -0.05%  '
$ws.Range('E25').Value = '  -2.23%  '
$cD = $ws.Range('D26')
$cD.NumberFormat = '@'
$cD.Value = '8.25'
$cD.Style = 'Normal'
$ws.Range('E26').Value = '  +5.03%  '
$cD = $ws.Range('D27')
$cD.NumberFormat = '@'
$cD.Value = '26.86'
$cD.Style = 'Normal'
$ws.Range('E27').Value = '  +2.77%  '
$ws.Range('E28').Value = '  +2.77%  '
$cD = $ws.Range('D29')
$cD.NumberFormat = '@'
$cD.Value = '7.23'
$cD.Style = 'Normal'
$ws.Range('E29').Value = '  -3.85%  '
$ws.Range('E30').Value = '  +0.05%  '
$ws.Range('E31').Value = '  -2.06%  '
$ws.Range('E32').Value = '  -0.85%  '
$cD = $ws.Range('D33')
$cD.NumberFormat = '@'
$cD.Value = '34.59'
$cD.Style = 'Normal'
$ws.Range('E33').Value = '  -0.60%  '
$ws.Range('E34').Value = '  +0.00%  '
$cD = $ws.Range('D35')
$cD.NumberFormat = '@'
$cD.Value = '50.38'
$cD.Style = 'Normal'
$ws.Range('E35').Value = '  -2.18%  '
$ws.Range('E36').Value = '  +1.05%  '
$ws.Range('E37').Value = '  +0.00%  '
$ws.Range('E38').Value = '  +2.28%  '
$cD = $ws.Range('D39')
$cD.NumberFormat = '@'
$cD.Value = '0.293'
$cD.Style = 'Normal'
$ws.Range('E39').Value = '  +7.21%  '
$cD = $ws.Range('D40')
$cD.NumberFormat = '@'
$cD.Value = '16.96'
$cD.Style = 'Normal'
$ws.Range('E40').Value = '  +0.77%  '
$cD = $ws.Range('D41')
$cD.NumberFormat = '@'
$cD.Value = '1.87'
$cD.Style = 'Normal'
$ws.Range('E41').Value = '  +1.09%  '
$ws.Range('E42').Value = '  -0.35%  '
$ws.Range('E43').Value = '  -1.23%  '
$cD = $ws.Range('D44')
$cD.NumberFormat = '@'
$cD.Value = '125.00'
$cD.Style = 'Normal'
$ws.Range('E44').Value = '  +0.11%  '
$ws.Range('E45').Value = '  +2.46%  '
$cD = $ws.Range('D46')
$cD.NumberFormat = '@'
$cD.Value = '21.83'
$cD.Style = 'Normal'
$ws.Range('E46').Value = '  +0.55%  '
$ws.Range('E47').Value = '  +3.62%  '
$cD = $ws.Range('D48')
$cD.NumberFormat = '@'
$cD.Value = '2.41'
$cD.Style = 'Normal'
$ws.Range('E48').Value = '  +1.55%  '
$cD = $ws.Range('D49')
$cD.NumberFormat = '@'
$cD.Value = '2.028.95'
$cD.Style = 'Normal'
$ws.Range('E49').Value = '  -0.45%  '
$cD = $ws.Range('D50')
$cD.NumberFormat = '@'
$cD.Value = '3.352.96'
$cD.Style = 'Normal'
$ws.Range('E50').Value = '  +2.38%  '
$cD = $ws.Range('D51')
$cD.NumberFormat = '@'
$cD.Value = '0.0319'
$cD.Style = 'Normal'
$ws.Range('E51').Value = '  -3.60%  '
